$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 10 for Testmail #8 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(10, 1).Value = "Kun je nagaan of dit nog leverbaar is?"
$logs.Cells.Item(10, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(10, 3).Value = "Testmail #8: Kun je nagaan of dit nog leverbaar is?"
$logs.Cells.Item(10, 4).Value = "Productinformatie"
$logs.Cells.Item(10, 5).Value = "Beste klant,
Bedankt voor uw e-mail. Om u beter van dienst te kunnen zijn, zou u ons kunnen voorzien van het specifieke product of de service waarover u informatie wenst over de leverbaarheid? Met deze informatie kunnen we gerichter voor u nagaan of het nog leverbaar is.
Met vriendelijke groet,
[Naam bedrijf]"
$logs.Cells.Item(10, 6).Value = "2025-07-29 21:44:05"
$logs.Cells.Item(10, 7).Value = "Ja"
$logs.Cells.Item(10, 8).Value = "Nee"
$logs.Cells.Item(10, 9).Value = "Ja"
$logs.Cells.Item(10, 10).Value = "Nee"

# --- Extend the conditional-formatting ranges from row 9 to row 10 ---
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "9")
    $newRange = $logs.Range($col + "2:" + $col + "10")
    $fc = $oldRange.FormatConditions.Item(1)
    $fc.ModifyAppliesToRange($newRange)
}

# --- Dashboard sheet: Productinformatie now outranks Intern verzoek (3 vs 2) ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(3, 1).Value = "Productinformatie"
$dash.Cells.Item(3, 2).Value = 3

$dash.Cells.Item(4, 1).Value = "Intern verzoek / Actie voor medewerker"
$dash.Cells.Item(4, 2).Value = 2
